$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 1.913996934890747
$ws.Range("D2").Value = 2
$ws.Range("E2").Value = 4745.685246913603
$ws.Range("F2").Value = 0.1085084843452944
$ws.Range("G2").Value = 0.1085084843452944
$ws.Range("H2").Value = 0.1085084843452944
$ws.Range("I2").Value = 0.1085084843452944
$ws.Range("J2").Value = 0.1085084843452944
$ws.Range("K2").Value = 0.1085084843452944
$ws.Range("L2").Value = 0.1085084843452944
$ws.Range("M2").Value = 0.1085084843452944
$ws.Range("N2").Value = 0.1085084843452944
$ws.Range("O2").Value = 0.1085084843452944
$ws.Range("P2").Value = 0.1085084843452944
$ws.Range("Q2").Value = 0.1085084843452944
$ws.Range("R2").Value = 0.1085084843452944
$ws.Range("S2").Value = 0.1085084843452944
$ws.Range("T2").Value = 0.1085084843452944
$ws.Range("U2").Value = 0.1085084843452944
$ws.Range("V2").Value = 0.1085084843452944
$ws.Range("W2").Value = 0.1085084843452944
$ws.Range("X2").Value = 0.1085084843452944
$ws.Range("Y2").Value = 0.1085084843452944
# Row 3
$ws.Range("C3").Value = 1.639001607894897
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = 4744.749434140682
$ws.Range("F3").Value = 0.10849024238091
$ws.Range("G3").Value = 0.10849024238091
$ws.Range("H3").Value = 0.10849024238091
$ws.Range("I3").Value = 0.10849024238091
$ws.Range("J3").Value = 0.10849024238091
$ws.Range("K3").Value = 0.10849024238091
$ws.Range("L3").Value = 0.10849024238091
$ws.Range("M3").Value = 0.10849024238091
$ws.Range("N3").Value = 0.10849024238091
$ws.Range("O3").Value = 0.10849024238091
$ws.Range("P3").Value = 0.10849024238091
$ws.Range("Q3").Value = 0.10849024238091
$ws.Range("R3").Value = 0.10849024238091
$ws.Range("S3").Value = 0.10849024238091
$ws.Range("T3").Value = 0.10849024238091
$ws.Range("U3").Value = 0.10849024238091
$ws.Range("V3").Value = 0.10849024238091
$ws.Range("W3").Value = 0.10849024238091
$ws.Range("X3").Value = 0.10849024238091
$ws.Range("Y3").Value = 0.10849024238091
# Row 4
$ws.Range("C4").Value = 1.743045091629028
$ws.Range("D4").Value = 2
$ws.Range("E4").Value = 4744.623808702443
$ws.Range("F4").Value = 0.1084877935419579
$ws.Range("G4").Value = 0.1084877935419579
$ws.Range("H4").Value = 0.1084877935419579
$ws.Range("I4").Value = 0.1084877935419579
$ws.Range("J4").Value = 0.1084877935419579
$ws.Range("K4").Value = 0.1084877935419579
$ws.Range("L4").Value = 0.1084877935419579
$ws.Range("M4").Value = 0.1084877935419579
$ws.Range("N4").Value = 0.1084877935419579
$ws.Range("O4").Value = 0.1084877935419579
$ws.Range("P4").Value = 0.1084877935419579
$ws.Range("Q4").Value = 0.1084877935419579
$ws.Range("R4").Value = 0.1084877935419579
$ws.Range("S4").Value = 0.1084877935419579
$ws.Range("T4").Value = 0.1084877935419579
$ws.Range("U4").Value = 0.1084877935419579
$ws.Range("V4").Value = 0.1084877935419579
$ws.Range("W4").Value = 0.1084877935419579
$ws.Range("X4").Value = 0.1084877935419579
$ws.Range("Y4").Value = 0.1084877935419579
# Row 5
$ws.Range("C5").Value = 1.781957149505615
$ws.Range("D5").Value = 2
$ws.Range("E5").Value = 4745.685246913603
$ws.Range("F5").Value = 0.1085084843452944
$ws.Range("G5").Value = 0.1085084843452944
$ws.Range("H5").Value = 0.1085084843452944
$ws.Range("I5").Value = 0.1085084843452944
$ws.Range("J5").Value = 0.1085084843452944
$ws.Range("K5").Value = 0.1085084843452944
$ws.Range("L5").Value = 0.1085084843452944
$ws.Range("M5").Value = 0.1085084843452944
$ws.Range("N5").Value = 0.1085084843452944
$ws.Range("O5").Value = 0.1085084843452944
$ws.Range("P5").Value = 0.1085084843452944
$ws.Range("Q5").Value = 0.1085084843452944
$ws.Range("R5").Value = 0.1085084843452944
$ws.Range("S5").Value = 0.1085084843452944
$ws.Range("T5").Value = 0.1085084843452944
$ws.Range("U5").Value = 0.1085084843452944
$ws.Range("V5").Value = 0.1085084843452944
$ws.Range("W5").Value = 0.1085084843452944
$ws.Range("X5").Value = 0.1085084843452944
$ws.Range("Y5").Value = 0.1085084843452944
# Row 6
$ws.Range("C6").Value = 1.54004693031311
$ws.Range("D6").Value = 2
$ws.Range("E6").Value = 4745.685246913603
$ws.Range("F6").Value = 0.1085084843452944
$ws.Range("G6").Value = 0.1085084843452944
$ws.Range("H6").Value = 0.1085084843452944
$ws.Range("I6").Value = 0.1085084843452944
$ws.Range("J6").Value = 0.1085084843452944
$ws.Range("K6").Value = 0.1085084843452944
$ws.Range("L6").Value = 0.1085084843452944
$ws.Range("M6").Value = 0.1085084843452944
$ws.Range("N6").Value = 0.1085084843452944
$ws.Range("O6").Value = 0.1085084843452944
$ws.Range("P6").Value = 0.1085084843452944
$ws.Range("Q6").Value = 0.1085084843452944
$ws.Range("R6").Value = 0.1085084843452944
$ws.Range("S6").Value = 0.1085084843452944
$ws.Range("T6").Value = 0.1085084843452944
$ws.Range("U6").Value = 0.1085084843452944
$ws.Range("V6").Value = 0.1085084843452944
$ws.Range("W6").Value = 0.1085084843452944
$ws.Range("X6").Value = 0.1085084843452944
$ws.Range("Y6").Value = 0.1085084843452944
# Row 7
$ws.Range("C7").Value = 1.682036399841309
$ws.Range("D7").Value = 2
$ws.Range("E7").Value = 4745.685246913603
$ws.Range("F7").Value = 0.1085084843452944
$ws.Range("G7").Value = 0.1085084843452944
$ws.Range("H7").Value = 0.1085084843452944
$ws.Range("I7").Value = 0.1085084843452944
$ws.Range("J7").Value = 0.1085084843452944
$ws.Range("K7").Value = 0.1085084843452944
$ws.Range("L7").Value = 0.1085084843452944
$ws.Range("M7").Value = 0.1085084843452944
$ws.Range("N7").Value = 0.1085084843452944
$ws.Range("O7").Value = 0.1085084843452944
$ws.Range("P7").Value = 0.1085084843452944
$ws.Range("Q7").Value = 0.1085084843452944
$ws.Range("R7").Value = 0.1085084843452944
$ws.Range("S7").Value = 0.1085084843452944
$ws.Range("T7").Value = 0.1085084843452944
$ws.Range("U7").Value = 0.1085084843452944
$ws.Range("V7").Value = 0.1085084843452944
$ws.Range("W7").Value = 0.1085084843452944
$ws.Range("X7").Value = 0.1085084843452944
$ws.Range("Y7").Value = 0.1085084843452944
# Row 8
$ws.Range("C8").Value = 1.742958545684814
$ws.Range("D8").Value = 2
$ws.Range("E8").Value = 4745.685246913603
$ws.Range("F8").Value = 0.1085084843452944
$ws.Range("G8").Value = 0.1085084843452944
$ws.Range("H8").Value = 0.1085084843452944
$ws.Range("I8").Value = 0.1085084843452944
$ws.Range("J8").Value = 0.1085084843452944
$ws.Range("K8").Value = 0.1085084843452944
$ws.Range("L8").Value = 0.1085084843452944
$ws.Range("M8").Value = 0.1085084843452944
$ws.Range("N8").Value = 0.1085084843452944
$ws.Range("O8").Value = 0.1085084843452944
$ws.Range("P8").Value = 0.1085084843452944
$ws.Range("Q8").Value = 0.1085084843452944
$ws.Range("R8").Value = 0.1085084843452944
$ws.Range("S8").Value = 0.1085084843452944
$ws.Range("T8").Value = 0.1085084843452944
$ws.Range("U8").Value = 0.1085084843452944
$ws.Range("V8").Value = 0.1085084843452944
$ws.Range("W8").Value = 0.1085084843452944
$ws.Range("X8").Value = 0.1085084843452944
$ws.Range("Y8").Value = 0.1085084843452944
# Row 9
$ws.Range("C9").Value = 1.659996747970581
$ws.Range("D9").Value = 2
$ws.Range("E9").Value = 4745.685246913603
$ws.Range("F9").Value = 0.1085084843452944
$ws.Range("G9").Value = 0.1085084843452944
$ws.Range("H9").Value = 0.1085084843452944
$ws.Range("I9").Value = 0.1085084843452944
$ws.Range("J9").Value = 0.1085084843452944
$ws.Range("K9").Value = 0.1085084843452944
$ws.Range("L9").Value = 0.1085084843452944
$ws.Range("M9").Value = 0.1085084843452944
$ws.Range("N9").Value = 0.1085084843452944
$ws.Range("O9").Value = 0.1085084843452944
$ws.Range("P9").Value = 0.1085084843452944
$ws.Range("Q9").Value = 0.1085084843452944
$ws.Range("R9").Value = 0.1085084843452944
$ws.Range("S9").Value = 0.1085084843452944
$ws.Range("T9").Value = 0.1085084843452944
$ws.Range("U9").Value = 0.1085084843452944
$ws.Range("V9").Value = 0.1085084843452944
$ws.Range("W9").Value = 0.1085084843452944
$ws.Range("X9").Value = 0.1085084843452944
$ws.Range("Y9").Value = 0.1085084843452944
# Row 10
$ws.Range("C10").Value = 1.651997327804565
$ws.Range("D10").Value = 2
$ws.Range("E10").Value = 4745.685246913603
$ws.Range("F10").Value = 0.1085084843452944
$ws.Range("G10").Value = 0.1085084843452944
$ws.Range("H10").Value = 0.1085084843452944
$ws.Range("I10").Value = 0.1085084843452944
$ws.Range("J10").Value = 0.1085084843452944
$ws.Range("K10").Value = 0.1085084843452944
$ws.Range("L10").Value = 0.1085084843452944
$ws.Range("M10").Value = 0.1085084843452944
$ws.Range("N10").Value = 0.1085084843452944
$ws.Range("O10").Value = 0.1085084843452944
$ws.Range("P10").Value = 0.1085084843452944
$ws.Range("Q10").Value = 0.1085084843452944
$ws.Range("R10").Value = 0.1085084843452944
$ws.Range("S10").Value = 0.1085084843452944
$ws.Range("T10").Value = 0.1085084843452944
$ws.Range("U10").Value = 0.1085084843452944
$ws.Range("V10").Value = 0.1085084843452944
$ws.Range("W10").Value = 0.1085084843452944
$ws.Range("X10").Value = 0.1085084843452944
$ws.Range("Y10").Value = 0.1085084843452944
# Row 11
$ws.Range("C11").Value = 1.836002111434937
$ws.Range("D11").Value = 2
$ws.Range("E11").Value = 4745.685246913603
$ws.Range("F11").Value = 0.1085084843452944
$ws.Range("G11").Value = 0.1085084843452944
$ws.Range("H11").Value = 0.1085084843452944
$ws.Range("I11").Value = 0.1085084843452944
$ws.Range("J11").Value = 0.1085084843452944
$ws.Range("K11").Value = 0.1085084843452944
$ws.Range("L11").Value = 0.1085084843452944
$ws.Range("M11").Value = 0.1085084843452944
$ws.Range("N11").Value = 0.1085084843452944
$ws.Range("O11").Value = 0.1085084843452944
$ws.Range("P11").Value = 0.1085084843452944
$ws.Range("Q11").Value = 0.1085084843452944
$ws.Range("R11").Value = 0.1085084843452944
$ws.Range("S11").Value = 0.1085084843452944
$ws.Range("T11").Value = 0.1085084843452944
$ws.Range("U11").Value = 0.1085084843452944
$ws.Range("V11").Value = 0.1085084843452944
$ws.Range("W11").Value = 0.1085084843452944
$ws.Range("X11").Value = 0.1085084843452944
$ws.Range("Y11").Value = 0.1085084843452944

Write-Host "Update complete: 230 cells updated"
